$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from column K (the last existing data column) to new column L,
# row by row (skip row 5, which gets no new cell in the target).
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("K6").Copy()
$ws.Range("L6").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("K7").Copy()
$ws.Range("L7").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("K8").Copy()
$ws.Range("L8").PasteSpecial(-4122)  # xlPasteFormats

# Add the new year column (2021) data in column L
$ws.Range("L4").Value = 2021
$ws.Range("L6").Value = 7.1
$ws.Range("L7").Value = 0.5
$ws.Range("L8").Value = "-"

# Update selection to match the target (activeCell N5)
$ws.Range("N5").Select()
